$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text used across rows 2-4 (shared string changes from "puri bhaji" to "pav bhaji")
$ws.Range("A2").Value = "pav bhaji"
$ws.Range("A3").Value = "pav bhaji"
$ws.Range("A4").Value = "pav bhaji"

# Update the amounts
$ws.Range("B2").Value = 50
$ws.Range("B3").Value = 60
$ws.Range("B4").Value = 45
